## update data dictionary and arcade code
##
## The "fad_choices" code list on the code_choices sheet gained a new entry
## for a hollow stem (code "HS"). Insert it right after the last existing
## fad_choices row ("R2" / "Armillaria root disease", row 47) and before the
## first "liana_load" row, pushing everything below down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("code_choices")

$ws.Rows.Item(47).Insert()
$ws.Range("A47").Value = "fad_choices"
$ws.Range("B47").Value = "HS"
$ws.Range("C47").Value = "Hollow stem"

# Leave the selection where the edit was made.
$ws.Range("B48").Select()
